# Trade #56 closed at 2026-02-18 00:21:16 - unknown UNKNOWN +0.000%
#
# 1) "Summary" sheet: refresh aggregate stats now that one more trade exists
#    and one more trade has closed.
# 2) "Strategy Status" sheet: refresh the "momentum" strategy row (row 11)
#    to reflect the extra trade/closed trade.
# 3) "All Trades" sheet: close trade #85 (row 86) and append the newly
#    opened trade #114 as a brand-new row (row 115).
# 4) "momentum" sheet (per-strategy log): same two edits as #3, mirrored
#    into this sheet's own row numbering (row 16 and new row 33).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.63
$summary.Range("B4").Value = 0.74
$summary.Range("B6").Value = 84
$summary.Range("B8").Value = 33
$summary.Range("B9").Value = 50

# ---------------------------------------------------------------------
# 2) Strategy Status - "momentum" row (row 11)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.7
$status.Range("D11").Value = 14
$status.Range("E11").Value = -0.3
$status.Range("F11").Value = -0.3
$status.Range("G11").Value = 14.29

# ---------------------------------------------------------------------
# 3) All Trades - close trade #85 (row 86), append trade #114 (row 115)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out the existing open trade in row 86.
$allTrades.Range("G86").Value = 0.84
$allTrades.Range("H86").Value = "CLOSED"
$allTrades.Range("I86").Value = -1.1765
$allTrades.Range("J86").Value = -0.01
$allTrades.Range("K86").Value = 99.7
$allTrades.Range("L86").Value = "early_exit"
$allTrades.Range("M86").Value = 0.12

# Append the newly-opened trade as row 115. Duplicate row 114 first so the
# unchanged columns (text dates/times, blanks, etc.) keep their original
# cell typing, then overwrite just the columns that actually differ.
$allTrades.Range("A114:Q114").Copy($allTrades.Range("A115:Q115"))
$allTrades.Range("A115").Value = 114
$allTrades.Range("C115").Value = "00:21:11"
$allTrades.Range("F115").Value = 0.85

# ---------------------------------------------------------------------
# 4) momentum sheet - mirror of the same two edits, own row numbering
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

# Close out trade #85, which lives in row 16 on this sheet.
$momentum.Range("G16").Value = 0.84
$momentum.Range("H16").Value = "CLOSED"
$momentum.Range("I16").Value = -1.1765
$momentum.Range("J16").Value = -0.01
$momentum.Range("K16").Value = 99.7
$momentum.Range("P16").Value = "early_exit"
$momentum.Range("Q16").Value = 0.12

# Append trade #114, which becomes row 33 on this sheet.
$momentum.Range("A32:Q32").Copy($momentum.Range("A33:Q33"))
$momentum.Range("A33").Value = 114
$momentum.Range("C33").Value = "00:21:11"
$momentum.Range("F33").Value = 0.85
